# "Sheet1" (the per-team pick/result helper tab) is removed from the
# workbook, leaving "Convoluted Scenario Table" as the sole remaining sheet.
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete() | Out-Null

# Re-activate the remaining sheet and restore its last active cell/selection.
$ws = $wb.Worksheets.Item("Convoluted Scenario Table")
$ws.Activate() | Out-Null
$ws.Range("D40").Select() | Out-Null
